# Apply weekly update to the "Achicoria" dataset.
# The edit reshuffles several rows' Fecha (D), Volumen (J), Precio minimo (K),
# Precio promedio ponderado (M), Origen (O) and Precio $/Kg (P) values so the
# sheet reflects the latest weekly data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44232

# Row 3
$ws.Range("D3").Value = 44188
$ws.Range("J3").Value = 210

# Row 4
$ws.Range("D4").Value = 44189
$ws.Range("J4").Value = 250

# Row 5
$ws.Range("D5").Value = 44231

# Row 6
$ws.Range("D6").Value = 44187
$ws.Range("J6").Value = 160

# Row 7
$ws.Range("D7").Value = 44215
$ws.Range("J7").Value = 250
$ws.Range("K7").Value = 5000
$ws.Range("M7").Value = 5500
$ws.Range("O7").Value = "Provincia de Quillota"
$ws.Range("P7").Value = 344

# Row 8
$ws.Range("D8").Value = 44292
$ws.Range("J8").Value = 90
$ws.Range("K8").Value = 6000
$ws.Range("M8").Value = 6000
$ws.Range("O8").Value = "Región Metropolitana"
$ws.Range("P8").Value = 375

# Row 9
$ws.Range("D9").Value = 44204
$ws.Range("J9").Value = 430

# Row 11
$ws.Range("D11").Value = 44210
$ws.Range("J11").Value = 340

# Row 12
$ws.Range("D12").Value = 44230

# Row 13
$ws.Range("D13").Value = 44186

# Row 14
$ws.Range("D14").Value = 44208
$ws.Range("J14").Value = 160

$wb.Save()
